$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("B14").Value = "-"
$ws.Range("B15").Value = "-"
$ws.Range("F18").Value = "-"
$ws.Range("F19").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("F21").Value = "-"
